$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking values: direct assignment ---
$ws.Range('D2').Value = '66.091.60'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '3.319.67'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E5').Value = '  +4.04%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -1.12%  '
$ws.Range('D9').Value = '3.313.59'
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('E11').Value = '  -0.53%  '
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('E14').Value = '  +1.05%  '
$ws.Range('D15').Value = '3.849.50'
$ws.Range('E15').Value = '  -0.24%  '
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').Value = '66.131.44'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('E18').Value = '  -1.35%  '
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('D20').Value = '3.321.74'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('E21').Value = '  -3.35%  '
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('E23').Value = '  +7.19%  '
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('E25').Value = '  -1.84%  '
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('E27').Value = '  +3.13%  '
$ws.Range('E28').Value = '  -1.27%  '
$ws.Range('E29').Value = '  +1.71%  '
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('E32').Value = '  +8.28%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('E34').Value = '  +9.04%  '
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('E38').Value = '  +2.06%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '3.709.69'
$ws.Range('E39').Value = '  -2.93%  '
$ws.Range('E40').Value = '  +6.65%  '
$ws.Range('E41').Value = '  +4.21%  '
$ws.Range('D42').Value = '0.0₃0714'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('E43').Value = '  -6.48%  '
$ws.Range('E44').Value = '  +5.17%  '
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('E46').Value = '  +4.30%  '
$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('E47').Value = '  +0.38%  '
$ws.Range('E48').Value = '  +1.77%  '
$ws.Range('E49').Value = '  -0.38%  '
$ws.Range('E50').Value = '  -1.18%  '
$ws.Range('E51').Value = '  +0.03%  '

# --- Numeric-looking values that must stay text (match source workbook's inlineStr type): ---
# set as a literal-text formula ( ="value" ), then immediately collapse the formula to a
# static value via Copy + PasteSpecial(xlPasteValues) on that same single cell, one at a time
# (multi-area ranges do not paste per-area reliably in this host, so this is done per cell).
$ws.Range('D4').Formula = '="1.00"'
$ws.Range('D4').Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('D5').Formula = '="188.52"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('D6').Formula = '="556.72"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('D7').Formula = '="1.00"'
$ws.Range('D7').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('D10').Formula = '="0.182"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('D11').Formula = '="0.584"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('D12').Formula = '="47.10"'
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('D14').Formula = '="8.66"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('D16').Formula = '="605.62"'
$ws.Range('D16').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('D18').Formula = '="17.96"'
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('D21').Formula = '="11.07"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('D22').Formula = '="0.905"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('D23').Formula = '="18.71"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('D24').Formula = '="5.06"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('D25').Formula = '="100.71"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('D26').Formula = '="3.96"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('D27').Formula = '="2.76"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('D28').Formula = '="5.92"'
$ws.Range('D28').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('D29').Formula = '="9.57"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('D30').Formula = '="8.72"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('D31').Formula = '="30.31"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('D32').Formula = '="6.79"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('D33').Formula = '="3.91"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('D34').Formula = '="573.11"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('D35').Formula = '="11.07"'
$ws.Range('D35').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('D37').Formula = '="1.00"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('D38').Formula = '="57.19"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('D40').Formula = '="34.10"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('D43').Formula = '="3.27"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('D44').Formula = '="3.36"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('D45').Formula = '="2.66"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('D46').Formula = '="3.42"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('D47').Formula = '="0.341"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('D48').Formula = '="0.0420"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('D49').Formula = '="0.129"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('D50').Formula = '="2.58"'
$ws.Range('D50').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('D51').Formula = '="1.00"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$excel.CutCopyMode = 0

Write-Host "Update complete"
